# Edit script for DateListNotibalityAnalyse.xlsx
# Shifts the date range forward and refreshes the data/labels, then
# removes the now-unused trailing rows (11-15), matching the new
# 9-row (A2:E10) data table used by the "start script for uwsgi" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1) Update the data table (rows 2-10) with the new dates/values
#    and the new Chinese labels in column E.
# ---------------------------------------------------------------

# Row 2
$ws.Range("A2").Value = 44593
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 8
$ws.Range("D2").Value = 8
$ws.Range("E2").Value = "家人"

# Row 3
$ws.Range("A3").Value = 44594
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 8
$ws.Range("D3").Value = 8
$ws.Range("E3").Value = "家人"

# Row 4 (also needs its A4 style fixed from the old wrap-text/date
# style to the plain date style used by the other rows)
$ws.Range("A2").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A4").Value = 44595
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 8
$ws.Range("D4").Value = 8
$ws.Range("E4").Value = "家人"

# Row 5 (E5 gains the wrap/vertical-center style used by the rest of
# column E)
$ws.Range("A5").Value = 44596
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 8
$ws.Range("D5").Value = 8
$ws.Range("E4").Copy()
$ws.Range("E5").PasteSpecial(-4122)
$ws.Range("E5").Value = "家人"

# Row 6
$ws.Range("A6").Value = 44597
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 8
$ws.Range("D6").Value = 8
$ws.Range("E6").Value = "家人"

# Row 7
$ws.Range("A7").Value = 44598
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 8
$ws.Range("D7").Value = 8
$ws.Range("E7").Copy()
$ws.Range("E7").PasteSpecial(-4122)
$ws.Range("E7").Value = "家人"

# Row 8 (E8 keeps its original style)
$ws.Range("A8").Value = 44599
$ws.Range("B8").Value = 3
$ws.Range("C8").Value = 8
$ws.Range("D8").Value = 6
$ws.Range("E8").Value = "预估"

# Row 9 (A9 style fix)
$ws.Range("A2").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("A9").Value = 44600
$ws.Range("B9").Value = 8
$ws.Range("C9").Value = 8
$ws.Range("D9").Value = 7
$ws.Range("E9").Value = "时间"

# Row 10 (A10 style fix)
$ws.Range("A2").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A10").Value = 44601
$ws.Range("B10").Value = 8
$ws.Range("C10").Value = 6
$ws.Range("D10").Value = 6
$ws.Range("E10").Value = "时间"

# ---------------------------------------------------------------
# 2) Remove the old trailing rows 11-15 (the table is now only
#    9 rows of data, A2:E10).
# ---------------------------------------------------------------
$ws.Range("A11:E15").Delete()

# ---------------------------------------------------------------
# 3) Point the line-chart series at the new, smaller A2:A10 /
#    B2:B10 / C2:C10 / D2:D10 ranges.
# ---------------------------------------------------------------
$co = $ws.ChartObjects(1)
$chart = $co.Chart
$s1 = $chart.SeriesCollection(1)
$s1.Formula = "=SERIES(Sheet1!`$B`$1,Sheet1!`$A`$2:`$A`$10,Sheet1!`$B`$2:`$B`$10,1)"
$s2 = $chart.SeriesCollection(2)
$s2.Formula = "=SERIES(Sheet1!`$C`$1,Sheet1!`$A`$2:`$A`$10,Sheet1!`$C`$2:`$C`$10,2)"
$s3 = $chart.SeriesCollection(3)
$s3.Formula = "=SERIES(Sheet1!`$D`$1,Sheet1!`$A`$2:`$A`$10,Sheet1!`$D`$2:`$D`$10,3)"

# ---------------------------------------------------------------
# 4) Match the saved selection/active cell (E16, just below the
#    new data table).
# ---------------------------------------------------------------
[void]$ws.Range("E16").Select()
